$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.746.91'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '2.099.25'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '227.24'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").Value = '61.83'
$ws.Range("E7").Value = '  +2.25%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("D10").Value = '0.0840'
$ws.Range("E10").Value = '  +0.28%  '

$ws.Range("E11").Value = '  -1.09%  '

$ws.Range("D12").Value = '15.77'
$ws.Range("E12").Value = '  +5.07%  '

$ws.Range("D13").Value = '2.411.29'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").Value = '0.800'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = '2.090.36'
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").Value = '38.716.92'
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("D19").Value = '71.56'
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("D20").Value = '6.01'
$ws.Range("E20").Value = '  -0.53%  '

$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").Value = '226.77'
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").Value = '2.33'
$ws.Range("E24").Value = '  -4.11%  '

$ws.Range("D25").Value = '2.32'
$ws.Range("E25").Value = '  -0.82%  '

$ws.Range("E26").Value = '  +1.59%  '

$ws.Range("D27").Value = '170.22'
$ws.Range("E27").Value = '  -0.26%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("D31").Value = '2.55'
$ws.Range("E31").Value = '  +8.52%  '

$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("E33").Value = '  +1.24%  '

$ws.Range("D34").Value = '4.79'
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("D35").Value = '7.16'

$ws.Range("D36").Value = '0.0612'
$ws.Range("E36").Value = '  -0.22%  '

$ws.Range("D37").Value = '2.36'
$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("D38").Value = '3.49'
$ws.Range("E38").Value = '  -1.49%  '

$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("D40").Value = '18.00'
$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("E41").Value = '  +3.07%  '

$ws.Range("D42").Value = '101.62'
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = '1.524.62'
$ws.Range("E43").Value = '  -1.33%  '

$ws.Range("E44").Value = '  +8.12%  '

$ws.Range("D45").Value = '2.80'
$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '7.74'
$ws.Range("E46").Value = '  +0.67%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.0910'
$ws.Range("E47").Value = '  -2.09%  '

$ws.Range("E48").Value = '  +4.33%  '

$ws.Range("E49").Value = '  +1.37%  '

$ws.Range("D50").Value = '2.95'
$ws.Range("E50").Value = '  -1.07%  '

$ws.Range("D51").Value = '2.298.73'
$ws.Range("E51").Value = '  -0.28%  '
